$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update Date_Created / Date_Expired timestamps (re-run offset ~4h27m later)
$ws1.Range("C2").Value2 = 43804.62040354496
$ws1.Range("D2").Value2 = 44170.62040354496
$ws1.Range("C3").Value2 = 43804.62042361854
$ws1.Range("D3").Value2 = 44170.62042361854
$ws1.Range("C4").Value2 = 43804.62044560831
$ws1.Range("D4").Value2 = 44170.62044560831
$ws1.Range("C5").Value2 = 43804.62046317705
$ws1.Range("D5").Value2 = 44170.62046317705
$ws1.Range("C6").Value2 = 43804.62048405896
$ws1.Range("D6").Value2 = 44170.62048405896
$ws1.Range("C7").Value2 = 43804.62050764143
$ws1.Range("D7").Value2 = 44170.62050764143
$ws1.Range("C8").Value2 = 43804.62053384448
$ws1.Range("D8").Value2 = 44170.62053384448
$ws1.Range("C9").Value2 = 43804.62055432207
$ws1.Range("D9").Value2 = 44170.62055432207
$ws1.Range("C10").Value2 = 43804.62057781933
$ws1.Range("D10").Value2 = 44170.62057781933
$ws1.Range("C11").Value2 = 43804.62059922087
$ws1.Range("D11").Value2 = 44170.62059922087
$ws1.Range("C12").Value2 = 43804.62062163772
$ws1.Range("D12").Value2 = 44170.62062163772
$ws1.Range("C13").Value2 = 43804.62064129587
$ws1.Range("D13").Value2 = 44170.62064129587
$ws1.Range("C14").Value2 = 43804.62067017692
$ws1.Range("D14").Value2 = 44170.62067017692
$ws1.Range("C15").Value2 = 43804.62068978882
$ws1.Range("D15").Value2 = 44170.62068978882
$ws1.Range("C16").Value2 = 43804.62071159387
$ws1.Range("D16").Value2 = 44170.62071159387
$ws1.Range("C17").Value2 = 43804.62073176025
$ws1.Range("D17").Value2 = 44170.62073176025
$ws1.Range("C18").Value2 = 43804.62076025998
$ws1.Range("D18").Value2 = 44170.62076025998
$ws1.Range("C19").Value2 = 43804.6207817073
$ws1.Range("D19").Value2 = 44170.6207817073
$ws1.Range("C20").Value2 = 43804.62080901489
$ws1.Range("D20").Value2 = 44170.62080901489
$ws1.Range("C21").Value2 = 43804.62082906438
$ws1.Range("D21").Value2 = 44170.62082906438

# Flagged sheet: append duplicate flagged entry from re-run
$wsFlagged = $wb.Worksheets.Item("Flagged")
$wsFlagged.Range("A3").Value = "Edwards"
$wsFlagged.Range("B3").Value = "Keith"
$wsFlagged.Range("C3").Value = "2019-12-05 14:53:53.686462"
$wsFlagged.Range("D3").Value = "2020-12-05 14:53:53.686462"
$wsFlagged.Range("E3").Value = "MS Center of NE New York"
$wsFlagged.Range("F3").Value = "Latham, NY"
